$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "Normal"
$ws.Range("D2").Value = "Sin modificaciones"

$ws.Range("C3").Value = "Comido"
$ws.Range("D3").Value = "asfvarfg"
